$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (A1:F7) holds parameter estimates as text strings (with
# specific formatting such as leading spaces / trailing zeros that must be
# preserved exactly), so force the target range to Text format before
# writing the new values - this keeps Excel from reinterpreting them as
# numbers.
$dataRng = $ws.Range("B2:F7")
$dataRng.NumberFormat = "@"

$ws.Range("B2").Value2 = "-0.114"
$ws.Range("C2").Value2 = "0.190"
$ws.Range("D2").Value2 = "-0.410"
$ws.Range("E2").Value2 = "-0.116"
$ws.Range("F2").Value2 = "0.187"

$ws.Range("B3").Value2 = " 1.027"
$ws.Range("C3").Value2 = "0.131"
$ws.Range("D3").Value2 = " 0.820"
$ws.Range("E3").Value2 = " 1.028"
$ws.Range("F3").Value2 = "1.232"

$ws.Range("B4").Value2 = "-0.058"
$ws.Range("C4").Value2 = "0.058"
$ws.Range("D4").Value2 = "-0.148"
$ws.Range("E4").Value2 = "-0.058"
$ws.Range("F4").Value2 = "0.033"

$ws.Range("B5").Value2 = " 0.073"
$ws.Range("C5").Value2 = "0.106"
$ws.Range("D5").Value2 = "-0.096"
$ws.Range("E5").Value2 = " 0.074"
$ws.Range("F5").Value2 = "0.238"

$ws.Range("B6").Value2 = "-0.154"
$ws.Range("C6").Value2 = "0.157"
$ws.Range("D6").Value2 = "-0.355"
$ws.Range("E6").Value2 = "-0.177"
$ws.Range("F6").Value2 = "0.159"

$ws.Range("B7").Value2 = " 0.094"
$ws.Range("C7").Value2 = "0.064"
$ws.Range("D7").Value2 = " 0.014"
$ws.Range("E7").Value2 = " 0.085"
$ws.Range("F7").Value2 = "0.206"

# Restore the original (default) cell style now that the text values are
# locked in, so no stray number-format style lingers on these cells.
$dataRng.Style = "Normal"
